$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (18) down through the new rows (19-28),
# mirroring the style-only "blank" D column cells and the numeric style on column C.
$ws.Range("A18:D18").Copy()
$ws.Range("A19:D28").PasteSpecial(-4122)

# New attendance/ranking entries
$fileNumbers = @(123473, 123474, 123475, 123476, 123477, 123478, 123479, 123480, 123481, 123483)
$names = @(
    "SCOTT, Michael",
    "BERTRAM, Nellie",
    "CALIFORNIA, Robert",
    "LEVINSON, Jan",
    "MILLER, Pete",
    "ANDERSON, Roy",
    "MINER, Charles",
    "BENNETT, Jo",
    "GREEN, Clark",
    "VANCE, Bob"
)
$points = @(7, 6, 0, 8, 6.5, 9, 4, 0, 6, 8)

$startRow = 19
for ($i = 0; $i -lt 10; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $fileNumbers[$i]
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $points[$i]
}

# Mirror the author's last interaction: selecting the full newly-added last row.
[void]$ws.Rows(28).Select()
